$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C) for rows 2-18 from 2023-09-21 to 2023-09-23
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45192
}
